$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pork Sandwich (row 2): mark Tested = Yes
$ws.Range("G2").Value = "Yes"

# Beef Sandwich (row 3): mark Recipe and Advancement = Yes
$ws.Range("E3").Value = "Yes"
$ws.Range("F3").Value = "Yes"

# Move active selection to E4
$ws.Range("E4").Select()
